$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Female, Czechia)
$ws.Range("D7").Value = 4803
$ws.Range("P7").Value = 145.4
$ws.Range("S7").Value = 3.2
$ws.Range("U7").Value = "145.4 (±152.5)"
$ws.Range("V7").Value = "3.1% (±3.2%)"
$ws.Range("X7").Value = 7.8
$ws.Range("Y7").Value = 8.300000000000001
$ws.Range("Z7").Value = "7.8(±8.3)"

# Row 38 (Male, Estonia)
$ws.Range("D38").Value = 1651
$ws.Range("P38").Value = 40.8
$ws.Range("R38").Value = 2.5
$ws.Range("U38").Value = "40.8 (±61.1)"
$ws.Range("V38").Value = "2.5% (±3.7%)"
$ws.Range("X38").Value = 19
$ws.Range("Z38").Value = "19.0(±28.5)"

# Row 43 (Male, Iceland)
$ws.Range("D43").Value = 158
$ws.Range("P43").Value = 9.4
$ws.Range("R43").Value = 6.3
$ws.Range("U43").Value = "9.4 (±14.9)"
$ws.Range("V43").Value = "6.3% (±9.7%)"
$ws.Range("X43").Value = 16.5
$ws.Range("Z43").Value = "16.5(±26.1)"

# Row 65 (Total, Czechia)
$ws.Range("D65").Value = 15106
$ws.Range("P65").Value = 603.6
$ws.Range("R65").Value = 4.2
$ws.Range("S65").Value = 3.3
$ws.Range("U65").Value = "603.6 (±469.6)"
$ws.Range("V65").Value = "4.2% (±3.3%)"
$ws.Range("X65").Value = 16.1
$ws.Range("Z65").Value = "16.1(±12.5)"

# Row 67 (Total, Estonia)
$ws.Range("D67").Value = 2353
$ws.Range("P67").Value = 116
$ws.Range("R67").Value = 5.2
$ws.Range("U67").Value = "116.0 (±85.7)"
$ws.Range("V67").Value = "5.2% (±3.9%)"
$ws.Range("X67").Value = 26.4
$ws.Range("Y67").Value = 19.4
$ws.Range("Z67").Value = "26.4(±19.4)"

# Row 72 (Total, Iceland)
$ws.Range("D72").Value = 269
$ws.Range("P72").Value = 16
$ws.Range("R72").Value = 6.3
$ws.Range("U72").Value = "16.0 (±15.0)"
$ws.Range("V72").Value = "6.3% (±5.9%)"
$ws.Range("X72").Value = 14.4
$ws.Range("Z72").Value = "14.4(±13.4)"
